$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.818.92'
$ws.Range('E2').Value = '  -0.62%  '

$ws.Range('D3').Value = '''1.735.51'

$ws.Range('D4').Value = '''0.9992'
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '''230.70'
$ws.Range('E5').Value = '  -2.70%  '

$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').Value = '''0.5209'
$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('D8').Value = '''0.2753'
$ws.Range('E8').Value = '  +1.64%  '

$ws.Range('D9').Value = '''39.30'
$ws.Range('E9').Value = '  -2.76%  '

$ws.Range('D10').Value = '''0.06124'
$ws.Range('E10').Value = '  -1.37%  '

$ws.Range('D11').Value = '''1.735.37'
$ws.Range('E11').Value = '  -1.26%  '

$ws.Range('D12').Value = '''0.07032'
$ws.Range('E12').Value = '  -0.03%  '

$ws.Range('D13').Value = '''14.95'
$ws.Range('E13').Value = '  -4.90%  '

$ws.Range('E14').Value = '  -3.48%  '

$ws.Range('D15').Value = '''4.515'
$ws.Range('E15').Value = '  +0.60%  '

$ws.Range('D16').Value = '''76.57'

$ws.Range('D17').Value = '''0.9988'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').Value = '''0.9997'
$ws.Range('E18').Value = '  +0.00%  '

$ws.Range('D19').Value = '''25.807.94'
$ws.Range('E19').Value = '  -0.69%  '

$ws.Range('D20').Value = '''11.44'
$ws.Range('E20').Value = '  -2.13%  '

$ws.Range('D21').Value = '''0.000006631'
$ws.Range('E21').Value = '  -1.04%  '

$ws.Range('D22').Value = '''1.959.84'
$ws.Range('E22').Value = '  -1.11%  '

$ws.Range('D23').Value = '''4.171'
$ws.Range('E23').Value = '  +2.10%  '

$ws.Range('D24').Value = '''8.743'
$ws.Range('E24').Value = '  +3.86%  '

$ws.Range('D25').Value = '''5.127'

$ws.Range('D26').Value = '''139.43'
$ws.Range('E26').Value = '  +1.70%  '

$ws.Range('D27').Value = '''1.501'
$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('E28').Value = '  -1.08%  '

$ws.Range('D29').Value = '''1.776'
$ws.Range('E29').Value = '  -2.68%  '

$ws.Range('D30').Value = '''101.84'
$ws.Range('E30').Value = '  -1.13%  '

$ws.Range('D31').Value = '''0.08261'
$ws.Range('E31').Value = '  -1.54%  '

$ws.Range('D32').Value = '''3.693'
$ws.Range('E32').Value = '  -0.14%  '

$ws.Range('D33').Value = '''3.496'
$ws.Range('E33').Value = '  +2.25%  '

$ws.Range('E34').Value = '  +1.01%  '

$ws.Range('D35').Value = '''2.602'
$ws.Range('E35').Value = '  -1.62%  '

$ws.Range('D36').Value = '''0.9718'
$ws.Range('E36').Value = '  -2.95%  '

$ws.Range('D37').Value = '''0.6136'
$ws.Range('E37').Value = '  +0.81%  '

$ws.Range('E38').Value = '  -1.80%  '

$ws.Range('D39').Value = '''0.01567'

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''1.909'
$ws.Range('E40').Value = '  -2.07%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''0.9994'
$ws.Range('E41').Value = '  -0.10%  '

$ws.Range('D42').Value = '''99.71'
$ws.Range('E42').Value = '  -2.92%  '

$ws.Range('D43').Value = '''0.3827'
$ws.Range('E43').Value = '  -1.45%  '

$ws.Range('D44').Value = '''0.7248'
$ws.Range('E44').Value = '  -3.59%  '

$ws.Range('D45').Value = '''4.987'
$ws.Range('E45').Value = '  +1.20%  '

$ws.Range('D46').Value = '''0.05342'
$ws.Range('E46').Value = '  -2.73%  '

$ws.Range('D47').Value = '''0.1127'
$ws.Range('E47').Value = '  +0.82%  '

$ws.Range('D48').Value = '''6.160'
$ws.Range('E48').Value = '  +0.95%  '

$ws.Range('D49').Value = '''53.04'
$ws.Range('E49').Value = '  +0.55%  '

$ws.Range('D50').Value = '''29.93'
$ws.Range('E50').Value = '  -0.88%  '

$ws.Range('D51').Value = '''7.605'
$ws.Range('E51').Value = '  +1.67%  '
